# Update Argent (Silver) prices and USD_CNY rate in the Solar_Prices workbook.
$wb = $excel.ActiveWorkbook

$updates = @{
    "Silver Rear_side"         = "5,455"
    "Silver Busbar front-side" = "8,167"
    "Silver finger front-side" = "8,217"
    "USD_CNY"                  = "7.2637"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range("B13")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$sheetName]
}
